$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.331632256507874
$ws.Range("B1").Value = 3.278643131256104
$ws.Range("C1").Value = 2.662660598754883
$ws.Range("D1").Value = 1.390368342399597
$ws.Range("E1").Value = 1.027259588241577
